$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D value (price text), only for rows whose D changed.
# These values must stay as TEXT (the sheet stores prices as inline strings,
# some of which look numeric, e.g. "553.91" or "20.00"), so we force the
# cell to Text format before assigning, then restore the default "Normal"
# style so no stray style index is left on the cell.
$dChanges = [ordered]@{
    2  = "60.204.40"
    3  = "2.421.29"
    5  = "553.91"
    6  = "137.11"
    10 = "5.68"
    14 = "2.852.30"
    15 = "60.109.82"
    17 = "2.428.85"
    20 = "328.12"
    21 = "6.76"
    23 = "65.38"
    24 = "0.177"
    25 = "8.69"
    30 = "170.45"
    31 = "6.11"
    33 = "0.405"
    34 = "18.55"
    39 = "327.97"
    41 = "144.96"
    43 = "20.00"
    50 = "4.66"
}

# Map of row -> new E value (volume text), for every changed row.
$eChanges = [ordered]@{
    2  = "  +0.16%  "
    3  = "  -0.05%  "
    4  = "  -0.03%  "
    5  = "  +0.05%  "
    6  = "  -1.16%  "
    7  = "  -0.01%  "
    8  = "  +1.66%  "
    9  = "  -1.17%  "
    10 = "  -1.80%  "
    12 = "  -1.47%  "
    13 = "  -0.19%  "
    14 = "  -0.08%  "
    15 = "  +0.15%  "
    16 = "  -0.98%  "
    17 = "  +0.21%  "
    18 = "  -0.73%  "
    19 = "  +2.17%  "
    20 = "  -1.34%  "
    21 = "  -0.41%  "
    22 = "  -0.10%  "
    23 = "  +0.34%  "
    24 = "  +3.41%  "
    25 = "  +1.08%  "
    26 = "  -0.02%  "
    27 = "  +2.68%  "
    28 = "  -2.04%  "
    29 = "  -1.03%  "
    30 = "  +0.37%  "
    31 = "  -3.18%  "
    32 = "  +1.30%  "
    33 = "  -4.13%  "
    34 = "  -0.89%  "
    36 = "  +1.27%  "
    37 = "  +0.06%  "
    38 = "  -0.04%  "
    39 = "  +2.38%  "
    40 = "  -0.90%  "
    41 = "  +3.65%  "
    42 = "  -1.08%  "
    43 = "  +2.36%  "
    44 = "  +0.63%  "
    45 = "  -1.21%  "
    46 = "  +0.03%  "
    47 = "  -1.29%  "
    49 = "  -2.66%  "
    50 = "  -0.55%  "
    51 = "  -0.46%  "
}

foreach ($row in $dChanges.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dChanges[$row]
    $cell.Style = "Normal"
}

foreach ($row in $eChanges.Keys) {
    $ws.Cells.Item($row, 5).Value = $eChanges[$row]
}
